$d = $word.ActiveDocument

# Update the title/date paragraph
$d.Content.Find.Execute("2023-12-14 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-15 Friday", 2) | Out-Null

# Update each table cell with its new expression
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "5+48="
$t.Cell(1, 2).Range.Text = "44-40="
$t.Cell(1, 3).Range.Text = "15-12="
$t.Cell(1, 4).Range.Text = "93+2="
$t.Cell(1, 5).Range.Text = "77-51="

$t.Cell(2, 1).Range.Text = "82-20="
$t.Cell(2, 2).Range.Text = "14+1="
$t.Cell(2, 3).Range.Text = "28+24="
$t.Cell(2, 4).Range.Text = "3+18="
$t.Cell(2, 5).Range.Text = "95-73="

$t.Cell(3, 1).Range.Text = "30+37="
$t.Cell(3, 2).Range.Text = "91-82="
$t.Cell(3, 3).Range.Text = "36+3="
$t.Cell(3, 4).Range.Text = "21+36="
$t.Cell(3, 5).Range.Text = "8+8="

$t.Cell(4, 1).Range.Text = "96-53="
$t.Cell(4, 2).Range.Text = "5+17="
$t.Cell(4, 3).Range.Text = "7-2="
$t.Cell(4, 4).Range.Text = "32+22="
$t.Cell(4, 5).Range.Text = "36+14="

$t.Cell(5, 1).Range.Text = "73-46="
$t.Cell(5, 2).Range.Text = "7+92="
$t.Cell(5, 3).Range.Text = "6+30="
$t.Cell(5, 4).Range.Text = "92-13="
$t.Cell(5, 5).Range.Text = "91-12="

$t.Cell(6, 1).Range.Text = "36+62="
$t.Cell(6, 2).Range.Text = "13+2="
$t.Cell(6, 3).Range.Text = "28+51="
$t.Cell(6, 4).Range.Text = "57-9="
$t.Cell(6, 5).Range.Text = "36+35="

$t.Cell(7, 1).Range.Text = "59-19="
$t.Cell(7, 2).Range.Text = "84-36="
$t.Cell(7, 3).Range.Text = "37-20="
$t.Cell(7, 4).Range.Text = "67+19="
$t.Cell(7, 5).Range.Text = "30+54="

$t.Cell(8, 1).Range.Text = "26+19="
$t.Cell(8, 2).Range.Text = "0+14="
$t.Cell(8, 3).Range.Text = "44-3="
$t.Cell(8, 4).Range.Text = "50-13="
$t.Cell(8, 5).Range.Text = "37+18="

$t.Cell(9, 1).Range.Text = "4+53="
$t.Cell(9, 2).Range.Text = "6+83="
$t.Cell(9, 3).Range.Text = "4+6="
$t.Cell(9, 4).Range.Text = "3+7="
$t.Cell(9, 5).Range.Text = "37+47="

$t.Cell(10, 1).Range.Text = "67-21="
$t.Cell(10, 2).Range.Text = "97-37="
$t.Cell(10, 3).Range.Text = "18+20="
$t.Cell(10, 4).Range.Text = "89-6="
$t.Cell(10, 5).Range.Text = "29+3="

$t.Cell(11, 1).Range.Text = "38-11="
$t.Cell(11, 2).Range.Text = "2+43="
$t.Cell(11, 3).Range.Text = "16+72="
$t.Cell(11, 4).Range.Text = "2+91="
$t.Cell(11, 5).Range.Text = "52-28="

$t.Cell(12, 1).Range.Text = "27+69="
$t.Cell(12, 2).Range.Text = "45-38="
$t.Cell(12, 3).Range.Text = "6+60="
$t.Cell(12, 4).Range.Text = "58-18="
$t.Cell(12, 5).Range.Text = "90-65="

$t.Cell(13, 1).Range.Text = "64+28="
$t.Cell(13, 2).Range.Text = "60+27="
$t.Cell(13, 3).Range.Text = "80-53="
$t.Cell(13, 4).Range.Text = "69-3="
$t.Cell(13, 5).Range.Text = "1+39="

$t.Cell(14, 1).Range.Text = "44+53="
$t.Cell(14, 2).Range.Text = "36-34="
$t.Cell(14, 3).Range.Text = "77+17="
$t.Cell(14, 4).Range.Text = "11+22="
$t.Cell(14, 5).Range.Text = "4+68="

$t.Cell(15, 1).Range.Text = "51-34="
$t.Cell(15, 2).Range.Text = "55-5="
$t.Cell(15, 3).Range.Text = "63+30="
$t.Cell(15, 4).Range.Text = "24+75="
$t.Cell(15, 5).Range.Text = "44+52="

$t.Cell(16, 1).Range.Text = "51-36="
$t.Cell(16, 2).Range.Text = "11+76="
$t.Cell(16, 3).Range.Text = "38+32="
$t.Cell(16, 4).Range.Text = "68-47="
$t.Cell(16, 5).Range.Text = "89-57="

$t.Cell(17, 1).Range.Text = "48+36="
$t.Cell(17, 2).Range.Text = "64-34="
$t.Cell(17, 3).Range.Text = "73-67="
$t.Cell(17, 4).Range.Text = "4+75="
$t.Cell(17, 5).Range.Text = "57+42="

$t.Cell(18, 1).Range.Text = "19+55="
$t.Cell(18, 2).Range.Text = "28+18="
$t.Cell(18, 3).Range.Text = "68-53="
$t.Cell(18, 4).Range.Text = "96-27="
$t.Cell(18, 5).Range.Text = "51-26="

$t.Cell(19, 1).Range.Text = "19+54="
$t.Cell(19, 2).Range.Text = "77-0="
$t.Cell(19, 3).Range.Text = "80-72="
$t.Cell(19, 4).Range.Text = "19-17="
$t.Cell(19, 5).Range.Text = "23+76="

$t.Cell(20, 1).Range.Text = "3+14="
$t.Cell(20, 2).Range.Text = "41+22="
$t.Cell(20, 3).Range.Text = "29+46="
$t.Cell(20, 4).Range.Text = "43-28="
$t.Cell(20, 5).Range.Text = "53+7="

